# Apply edit: insert two new rows of data into Sheet1 so that the table
# gains "un_franzosa_ControlvsCD_ConvCD" (right after ControlvsCD_Age /
# before ControlvsCD_Fp) and "un_franzosa_ControlvsUC_ConvUC" (right
# after ControlvsUC_Age / before ControlvsUC_Fp), shifting the rest of
# the table down accordingly. This grows the used range from H26 to H28.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Insert first new row: "un_franzosa_ControlvsCD_ConvCD" ---
# Currently row 9 holds "un_franzosa_ControlvsCD_Fp"; push it (and everything
# below) down by inserting a new blank row at position 9.
$ws.Rows.Item(9).Insert() | Out-Null

$ws.Cells.Item(9, 1).Value = "un_franzosa_ControlvsCD_ConvCD"
$ws.Cells.Item(9, 2).Value = 0
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(9, 4).Value = 0
$ws.Cells.Item(9, 5).Value = 0.45
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.55
$ws.Cells.Item(9, 8).Value = 0.55

# --- Insert second new row: "un_franzosa_ControlvsUC_ConvUC" ---
# After the first insertion, "un_franzosa_ControlvsUC_Age" is now at row 14
# and "un_franzosa_ControlvsUC_Fp" is now at row 15. Insert a new blank row
# at position 15 (before Fp) to hold the new ConvUC entry.
$ws.Rows.Item(15).Insert() | Out-Null

$ws.Cells.Item(15, 1).Value = "un_franzosa_ControlvsUC_ConvUC"
$ws.Cells.Item(15, 2).Value = 0
$ws.Cells.Item(15, 3).Value = 0
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 0.3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 0.7
$ws.Cells.Item(15, 8).Value = 0.7
